$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.5144553327538404
$ws.Range("C2").Value = 0.143483120846156
$ws.Range("E2").Value = 0.6295212807565846
$ws.Range("F2").Value = 1.787587014989356
$ws.Range("G2").Value = 0.2331284652465868
$ws.Range("H2").Value = 0.4045847091509742
$ws.Range("I2").Value = 0.2862892952744467
$ws.Range("J2").Value = 0.0215090097212709
$ws.Range("M2").Value = 0.5033101217576785
$ws.Range("N2").Value = 0.9762013637867071
$ws.Range("O2").Value = 1.190126071975598
$ws.Range("B3").Value = 0.4489467493570487
$ws.Range("C3").Value = 0.1295938597751274
$ws.Range("E3").Value = 0.6248243034619492
$ws.Range("F3").Value = 1.77476356758109
$ws.Range("G3").Value = 0.2315883309147395
$ws.Range("H3").Value = 0.4071382128668262
$ws.Range("I3").Value = 0.2904288655679181
$ws.Range("J3").Value = 0.0215527542838263
$ws.Range("M3").Value = 0.4722958032991542
$ws.Range("N3").Value = 0.9784584358061963
$ws.Range("O3").Value = 1.191976016832186
$ws.Range("B4").Value = 0.408576771306997
$ws.Range("C4").Value = 0.1210138013745592
$ws.Range("E4").Value = 0.622233778496863
$ws.Range("F4").Value = 1.767937267232014
$ws.Range("G4").Value = 0.2308625341639612
$ws.Range("H4").Value = 0.408917673721426
$ws.Range("I4").Value = 0.2931774368730871
$ws.Range("J4").Value = 0.02158559926387937
$ws.Range("M4").Value = 0.4533953736997844
$ws.Range("N4").Value = 0.9802148788764669
$ws.Range("O4").Value = 1.19401285737618
$ws.Range("B5").Value = 0.3920898038851135
$ws.Range("C5").Value = 0.1175045224172635
$ws.Range("E5").Value = 0.6212519556594884
$ws.Range("F5").Value = 1.765418832972941
$ws.Range("G5").Value = 0.2306219167950161
$ws.Range("H5").Value = 0.4096960353688672
$ws.Range("I5").Value = 0.2943494442234424
$ws.Range("J5").Value = 0.02160049399746811
$ws.Range("M5").Value = 0.4457295734653428
$ws.Range("N5").Value = 0.9810240397697214
$ws.Range("O5").Value = 1.19506917257813
$ws.Range("B6").Value = 0.3893500242425034
$ws.Range("C6").Value = 0.1169210411391504
$ws.Range("E6").Value = 0.6210933857890453
$ws.Range("F6").Value = 1.765016553183628
$ws.Range("G6").Value = 0.2305852895297207
$ws.Range("H6").Value = 0.4098284961527341
$ws.Range("I6").Value = 0.294547190161639
$ws.Range("J6").Value = 0.02160305861194267
$ws.Range("M6").Value = 0.4444588775147977
$ws.Range("N6").Value = 0.9811640473988632
$ws.Range("O6").Value = 1.195258232140304
$ws.Range("B7").Value = 0.4083545659554204
$ws.Range("C7").Value = 0.1209665256643149
$ws.Range("E7").Value = 0.6222202382614697
$ws.Range("F7").Value = 1.767902236534837
$ws.Range("G7").Value = 0.2308590659878433
$ws.Range("H7").Value = 0.4089279554946472
$ws.Range("I7").Value = 0.2931930327602998
$ws.Range("J7").Value = 0.02158579401870142
$ws.Range("M7").Value = 0.453291842552936
$ws.Range("N7").Value = 0.9802254130682755
$ws.Range("O7").Value = 1.194026187347305
$ws.Range("B8").Value = 0.4918992760025276
$ws.Range("C8").Value = 0.1387050347939578
$ws.Range("E8").Value = 0.6278408866844103
$ws.Range("F8").Value = 1.78294808863447
$ws.Range("G8").Value = 0.2325517156295902
$ws.Range("H8").Value = 0.4054212445880623
$ws.Range("I8").Value = 0.287673641494143
$ws.Range("J8").Value = 0.02152285350386762
$ws.Range("M8").Value = 0.4925870704318811
$ws.Range("N8").Value = 0.9769028079333779
$ws.Range("O8").Value = 1.190576772166125
$ws.Range("B9").Value = 0.6545150870073826
$ws.Range("C9").Value = 0.1730686826175827
$ws.Range("E9").Value = 0.6411903915458481
$ws.Range("F9").Value = 1.820768502141576
$ws.Range("G9").Value = 0.2376219936468331
$ws.Range("H9").Value = 0.4002233684093355
$ws.Range("I9").Value = 0.2784947202575481
$ws.Range("J9").Value = 0.02144670930521286
$ws.Range("M9").Value = 0.5707598490532746
$ws.Range("N9").Value = 0.9733193365094337
$ws.Range("O9").Value = 1.190975070110298
$ws.Range("B10").Value = 0.7731985689076737
$ws.Range("C10").Value = 0.1980489613613088
$ws.Range("E10").Value = 0.6524179189058188
$ws.Range("F10").Value = 1.853637867193413
$ws.Range("G10").Value = 0.2424245055813827
$ws.Range("H10").Value = 0.3974279859112642
$ws.Range("I10").Value = 0.2727579498780841
$ws.Range("J10").Value = 0.0214193183045559
$ws.Range("M10").Value = 0.6288581517382852
$ws.Range("N10").Value = 0.972464043256366
$ws.Range("O10").Value = 1.195655748490964
$ws.Range("B11").Value = 0.8270093143064514
$ws.Range("C11").Value = 0.2093533064017379
$ws.Range("E11").Value = 0.6578341763028916
$ws.Range("F11").Value = 1.869697893469521
$ws.Range("G11").Value = 0.2448454034763046
$ws.Range("H11").Value = 0.3963785785008298
$ws.Range("I11").Value = 0.2703677413222216
$ws.Range("J11").Value = 0.02141300045667194
$ws.Range("M11").Value = 0.655430116526631
$ws.Range("N11").Value = 0.9724589342257417
$ws.Range("O11").Value = 1.198742692573177
$ws.Range("B12").Value = 0.8473592316366307
$ws.Range("C12").Value = 0.2136252256700857
$ws.Range("E12").Value = 0.6599295624169059
$ws.Range("F12").Value = 1.875938817547024
$ws.Range("G12").Value = 0.245796257113355
$ws.Range("H12").Value = 0.3960131512319549
$ws.Range("I12").Value = 0.2694942729758445
$ws.Range("J12").Value = 0.02141148644921742
$ws.Range("M12").Value = 0.6655123826936915
$ws.Range("N12").Value = 0.9725120414057926
$ws.Range("O12").Value = 1.200049687346706
$ws.Range("B13").Value = 0.8429777364382858
$ws.Range("C13").Value = 0.2127055865662726
$ws.Range("E13").Value = 0.6594763108931829
$ws.Range("F13").Value = 1.874587635046282
$ws.Range("G13").Value = 0.2455899543234779
$ws.Range("H13").Value = 0.3960904310048221
$ws.Range("I13").Value = 0.2696809808435603
$ws.Range("J13").Value = 0.02141177352634038
$ws.Range("M13").Value = 0.6633401027789603
$ws.Range("N13").Value = 0.9724981585221286
$ws.Range("O13").Value = 1.199762058337313
$ws.Range("B14").Value = 0.8286840637442765
$ws.Range("C14").Value = 0.2097049374535516
$ws.Range("E14").Value = 0.6580056760154136
$ws.Range("F14").Value = 1.870208144223099
$ws.Range("G14").Value = 0.2449229462915241
$ws.Range("H14").Value = 0.3963478740156603
$ws.Range("I14").Value = 0.2702952456860821
$ws.Range("J14").Value = 0.02141285832123252
$ws.Range("M14").Value = 0.6562591918017375
$ws.Range("N14").Value = 0.9724622012366382
$ws.Range("O14").Value = 1.198847451533197
$ws.Range("B15").Value = 0.8199252157563137
$ws.Range("C15").Value = 0.2078658006250294
$ws.Range("E15").Value = 0.6571106464593086
$ws.Range("F15").Value = 1.867546332582307
$ws.Range("G15").Value = 0.244518831310188
$ws.Range("H15").Value = 0.3965097276377776
$ws.Range("I15").Value = 0.2706756253329807
$ws.Range("J15").Value = 0.02141363704210342
$ws.Range("M15").Value = 0.6519245246520882
$ws.Range("N15").Value = 0.9724473393934829
$ws.Range("O15").Value = 1.198305214355742
$ws.Range("B16").Value = 0.7696781969265771
$ws.Range("C16").Value = 0.1973089799409422
$ws.Range("E16").Value = 0.6520701643229287
$ws.Range("F16").Value = 1.852610600290674
$ws.Range("G16").Value = 0.2422710598300739
$ws.Range("H16").Value = 0.3975010393140508
$ws.Range("I16").Value = 0.2729185798749079
$ws.Range("J16").Value = 0.02141985437768668
$ws.Range("M16").Value = 0.6271244422131446
$ws.Range("N16").Value = 0.9724720887408012
$ws.Range("O16").Value = 1.195473308067932
$ws.Range("B17").Value = 0.7388065085892777
$ws.Range("C17").Value = 0.1908173350853701
$ws.Range("E17").Value = 0.6490570555371633
$ws.Range("F17").Value = 1.843731742880649
$ws.Range("G17").Value = 0.240952727942755
$ws.Range("H17").Value = 0.3981660977183452
$ws.Range("I17").Value = 0.2743508382061464
$ws.Range("J17").Value = 0.02142523828502263
$ws.Range("M17").Value = 0.611946610589186
$ws.Range("N17").Value = 0.9725854911995384
$ws.Range("O17").Value = 1.193981540728601
$ws.Range("B18").Value = 0.7210331629316897
$ws.Range("C18").Value = 0.1870779458505751
$ws.Range("E18").Value = 0.6473530629189383
$ws.Range("F18").Value = 1.838729106122372
$ws.Range("G18").Value = 0.2402166830152908
$ws.Range("H18").Value = 0.3985695384164174
$ws.Range("I18").Value = 0.2751952881057669
$ws.Range("J18").Value = 0.02142891341240549
$ws.Range("M18").Value = 0.6032301892858527
$ws.Range("N18").Value = 0.972686868223775
$ws.Range("O18").Value = 1.193213637864517
$ws.Range("B19").Value = 0.7150125736605446
$ws.Range("C19").Value = 0.1858109048058907
$ws.Range("E19").Value = 0.6467811138910449
$ws.Range("F19").Value = 1.837053201647834
$ws.Range("G19").Value = 0.2399712835172494
$ws.Range("H19").Value = 0.3987097288237322
$ws.Range("I19").Value = 0.2754847486832901
$ws.Range("J19").Value = 0.02143025724471492
$ws.Range("M19").Value = 0.6002812864803388
$ws.Range("N19").Value = 0.9727274070190788
$ws.Range("O19").Value = 1.192969107852065
$ws.Range("B20").Value = 0.7420945973327662
$ws.Range("C20").Value = 0.1915089598074076
$ws.Range("E20").Value = 0.6493747979714968
$ws.Range("F20").Value = 1.844666122696893
$ws.Range("G20").Value = 0.2410907653927694
$ws.Range("H20").Value = 0.3980931362820996
$ws.Range("I20").Value = 0.2741962337133153
$ws.Range("J20").Value = 0.02142460532692958
$ws.Range("M20").Value = 0.6135609265319459
$ws.Range("N20").Value = 0.9725696790890908
$ws.Range("O20").Value = 1.194131012244327
$ws.Range("B21").Value = 0.8328832062860556
$ws.Range("C21").Value = 0.2105865415038295
$ws.Range("E21").Value = 0.6584364333212847
$ws.Range("F21").Value = 1.871490181818871
$ws.Range("G21").Value = 0.245117935773834
$ws.Range("H21").Value = 0.3962713893253351
$ws.Range("I21").Value = 0.2701139614670751
$ws.Range("J21").Value = 0.02141251588938786
$ws.Range("M21").Value = 0.6583384870165361
$ws.Range("N21").Value = 0.9724712702950598
$ws.Range("O21").Value = 1.199112344904904
$ws.Range("B22").Value = 0.8920605714012595
$ws.Range("C22").Value = 0.2230034828431258
$ws.Range("E22").Value = 0.6646173231988683
$ws.Range("F22").Value = 1.889950046524774
$ws.Range("G22").Value = 0.2479488226008471
$ws.Range("H22").Value = 0.3952670650173644
$ws.Range("I22").Value = 0.2676304933056528
$ws.Range("J22").Value = 0.02140973262523183
$ws.Range("M22").Value = 0.6877197919563116
$ws.Range("N22").Value = 0.97272770351654
$ws.Range("O22").Value = 1.20317266120847
$ws.Range("B23").Value = 0.8604914085728979
$ws.Range("C23").Value = 0.2163811136076674
$ws.Range("E23").Value = 0.6612948176483897
$ws.Range("F23").Value = 1.880012661723882
$ws.Range("G23").Value = 0.2464196781641732
$ws.Range("H23").Value = 0.3957860444767789
$ws.Range("I23").Value = 0.2689390522017732
$ws.Range("J23").Value = 0.02141075142007587
$ws.Range("M23").Value = 0.6720279307377552
$ws.Range("N23").Value = 0.9725615475729228
$ws.Range("O23").Value = 1.200931855459714
$ws.Range("B24").Value = 0.740608129563725
$ws.Range("C24").Value = 0.1911962988534128
$ws.Range("E24").Value = 0.6492310585116385
$ws.Range("F24").Value = 1.844243372060745
$ws.Range("G24").Value = 0.241028290513313
$ws.Range("H24").Value = 0.3981260564538331
$ws.Range("I24").Value = 0.2742660649183613
$ws.Range("J24").Value = 0.02142488968099165
$ws.Range("M24").Value = 0.6128310646152215
$ws.Range("N24").Value = 0.9725767150307121
$ws.Range("O24").Value = 1.194063156680983
$ws.Range("B25").Value = 0.610658425533245
$ws.Range("C25").Value = 0.1638185393583456
$ws.Range("E25").Value = 0.6373297497017347
$ws.Range("F25").Value = 1.809645461075405
$ws.Range("G25").Value = 0.2360618628445792
$ws.Range("H25").Value = 0.4014498000757172
$ws.Range("I25").Value = 0.2808013415736355
$ws.Range("J25").Value = 0.02146227465209272
$ws.Range("M25").Value = 0.5494942228231139
$ws.Range("N25").Value = 0.9739759079397174
$ws.Range("O25").Value = 1.190098279898393
